# Fixed mistake with "actual" sim pars
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")
$ws.Activate()

# Row 21 ("actual" scenario, last group of replicate columns T:V): 40 -> 50
$ws.Range("T21:V21").Value = 50

# Row 22 ("actual" scenario, last group of replicate columns T:V): 60 -> 70
$ws.Range("T22:V22").Value = 70

# Row 30 held the id labels for the (now removed) extra "Actual ... 2023-10-24" columns.
# Clear them out entirely so the shared-string entries are dropped too.
$ws.Range("T30:V30").ClearContents()

# Restore selection to the last-edited cell, as in the saved file.
$ws.Range("V21").Select()
